$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.566.76"
$ws.Range("E2").Value = "  +3.57%  "
$ws.Range("D3").Value = "1.695.08"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3953"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4017"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("E9").Value = "  +7.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08762"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001317"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.600"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").Value = "1.692.92"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.876"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "24.566.05"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.067"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.331"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.188"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.22%  "
$ws.Range("D32").Value = "1.879.74"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.094"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.358"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08532"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.970"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2734"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09012"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.463"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7676"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7194"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.226"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08032"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.39%  "
